# Applies the diff between the original mclone-sw-structure.pptx and the
# edited revision: resize several "CustomShape" boxes by a hair (re-export
# rounding drift) and fix up a few labels in the "Util" and "Data" boxes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) CustomShape 1 "Util" -------------------------------------------------
$sh1 = $s.Shapes.Item(1)
$sh1.Width  = 226.1764
$sh1.Height = 421.6252

$tr1 = $sh1.TextFrame.TextRange
$run = $tr1.Paragraphs(2,1)
$run.Text = "XplaceholderX"
$run = $tr1.Paragraphs(2,1)
$run.Text = "distPoints"

# ---------------------------------------------------------------------
# 2) CustomShape 2 "Simulation" ------------------------------------------
$sh2 = $s.Shapes.Item(2)
$sh2.Width  = 189.4111
$sh2.Height = 230.1449

# ---------------------------------------------------------------------
# 3) CustomShape 3 "Control" ----------------------------------------------
$sh3 = $s.Shapes.Item(3)
$sh3.Width  = 208.8000
$sh3.Height = 249.1087

# ---------------------------------------------------------------------
# 4) CustomShape 4 "Data" --------------------------------------------------
$sh4 = $s.Shapes.Item(4)
$sh4.Width  = 225.8079
$sh4.Height = 322.4977

$tr4 = $sh4.TextFrame.TextRange

$run = $tr4.Paragraphs(3,1)
$run.Text = "XplaceholderX"
$run = $tr4.Paragraphs(3,1)
$run.Text = "cellsList"

$run = $tr4.Paragraphs(6,1)
$run.Text = "XplaceholderX"
$run = $tr4.Paragraphs(6,1)
$run.Text = "Point3D"
$run.Font.Color.RGB = 0

$run = $tr4.Paragraphs(7,1)
$run.Text = "XplaceholderX"
$run = $tr4.Paragraphs(7,1)
$run.Text = "Matrix4"
$run.Font.Color.RGB = 0

$run = $tr4.Paragraphs(9,1)
$run.Font.Color.RGB = 0

# Drop the now-superfluous trailing empty paragraph after "fileManager".
$tr4.Paragraphs(10,1).Delete()

# ---------------------------------------------------------------------
# 5) CustomShape 5 "Viewer" -------------------------------------------------
$sh5 = $s.Shapes.Item(5)
$sh5.Width  = 189.4111
$sh5.Height = 62.1922

# ---------------------------------------------------------------------
# 6) CustomShape 6 "Libs" ---------------------------------------------------
$sh6 = $s.Shapes.Item(6)
$sh6.Width  = 189.4111
$sh6.Height = 129.3733
